$d = $word.ActiveDocument

function Get-ParaRangeByText($anchorText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($anchorText)
    if (-not $ok) {
        throw "Anchor text not found: $anchorText"
    }
    $para = $rng.Paragraphs(1)
    return $para.Range
}

function Replace-ParaXML($anchorText, $innerXml) {
    $pRange = Get-ParaRangeByText $anchorText
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pRange.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) Paragraphs whose runs were merely split across multiple <w:r>
#    get collapsed back into a single run. Plain Find/Replace of the
#    whole-paragraph text (matching the final single-run content)
#    naturally produces a single merged run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "What would you choose as the key measure of success of this experiment in encouraging driver partners to serve both cities, and why would you choose this metric?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "What would you choose as the key measure of success of this experiment in encouraging driver partners to serve both cities, and why would you choose this metric?",
    2) | Out-Null

$d.Content.Find.Execute(
    "Describe a practical experiment you would design to compare the effectiveness of the proposed change in relation to the key measure of success. Please provide details on:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Describe a practical experiment you would design to compare the effectiveness of the proposed change in relation to the key measure of success. Please provide details on:",
    2) | Out-Null

$d.Content.Find.Execute(
    "what statistical test(s) you will conduct to verify the significance of the observation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "what statistical test(s) you will conduct to verify the significance of the observation",
    2) | Out-Null

$d.Content.Find.Execute(
    "how you would interpret the results and provide recommendations to the city operations team along with any caveats.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "how you would interpret the results and provide recommendations to the city operations team along with any caveats.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) "Toll Amnesty Week" paragraph: trim the trailing sentence, add a
#    blank spacer paragraph, then a new paragraph describing the ratio
#    methodology (several runs).
# ---------------------------------------------------------------------
$tollInner = (
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">One practical method would be a Toll Amnesty Week.  This would be well publicized to drivers in advance of the experiment and would involve reimbursement of tolls during weekdays for a routine week with no holidays, etc.  </w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">Since the goal is to encourage driver partners to serve both cities, we would calculate a ratio for drivers.  </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>Using trip origination data</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> during </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>a</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> period</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">, we can </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">calculate </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>ratios of trip originations</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> in each city for each day.</w:t></w:r>' +
    '</w:p>'
)
Replace-ParaXML "One practical method would be a Toll Amnesty Week" $tollInner

# ---------------------------------------------------------------------
# 3) "Th" + "is" + ... paragraph -> replaced with null-hypothesis text,
#    followed by a brand-new paragraph about the control comparison.
# ---------------------------------------------------------------------
$hypoInner = (
    '<w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>The null hypothesis would be that driver' + [char]0x2019 + 's ratios of trip origination cities remains the same.  The alternative hypothesis would be that the ratios decrease (flatten).</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>The prior week or weeks would be a valid control comparison.</w:t></w:r>' +
    '</w:p>'
)
Replace-ParaXML "would be company managed and promoted with established measurement protocols" $hypoInner

# ---------------------------------------------------------------------
# 4) "This is basic hypothesis testing whereby ..." paragraph rewritten
#    with the permutation-test explanation (several runs).
# ---------------------------------------------------------------------
$permInner = (
    '<w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>Using a permutation test we</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> calculate the likelihood that </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>improved</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> results were possible considering our null hypothesis that there was no effect.</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">  If we have greater than 95% confidence that the results were not random, (p-value &lt; .05) then we have significance.</w:t></w:r>' +
    '</w:p>'
)
Replace-ParaXML "This is basic hypothesis testing whereby" $permInner

# ---------------------------------------------------------------------
# 5) "If we set our significance level ..." and the following "We would
#    present ..." paragraph are merged into a single paragraph with the
#    new significance-measures wording.
# ---------------------------------------------------------------------
$presentInner = (
    '<w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">We would present </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>significance measures</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">a summary of </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>any financial impact and customer satisfaction impact observed.</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">  These measures could also have permutation tests applied.</w:t></w:r>' +
    '</w:p>'
)

# Remove the old "We would present ..." paragraph entirely (its content
# is superseded / merged into the paragraph above), then rewrite the
# "If we set ..." paragraph.
$oldPresentRange = Get-ParaRangeByText "a financial summary of the impact of making a similar change on an ongoing basis."
$oldPresentRange.Text = ""
$oldPresentRange.Delete()

Replace-ParaXML "If we set our significance level at a typical .05" $presentInner

Write-Output "Done"
